$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 ("Measures reported, House joint resolutions"): the Senate count was
# mistakenly entered as the text "l" (lowercase L) instead of the number 11.
$ws.Range("B25").Value = 11

# Row 42: correct the misspelled label "Vetoes overriden" to "Vetoes overridden".
$ws.Range("A42").Value = "Vetoes overridden  "
